# Updated symbol list (coin prices / names / links) as described by the diff.
#
# Every touched cell in this sheet is stored as TEXT (t="inlineStr" in the
# original file, round-tripped to a shared string by this engine) even
# though many of the "Price" column values look like plain numbers
# ("250.61", "0.006771", ...). A plain `Range.Value = "250.61"` assignment
# lets Excel's usual type-inference kick in and silently turns the cell
# into a numeric cell (t="n", and it also drifts to a binary float like
# 250.61000000000001). To keep those cells as text - matching the target
# workbook exactly - we temporarily force the cell to the "Text" number
# format ("@") before writing the value, then call ClearFormats() to drop
# the temporary format again so no stray style index is left behind on
# the cell (the original cells carry no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values, keyed by A1 reference, taken from the diff.
$newValues = [ordered]@{
    "D2"  = "250.61"
    "D3"  = "22.90"
    "D4"  = "5.423"
    "D5"  = "0.05668"
    "D6"  = "3.414"
    "D7"  = "6.380"
    "D8"  = "0.8135"
    "D9"  = "0.9248"
    "D10" = "0.1441"
    "D11" = "0.07442"
    "D12" = "0.03112"
    "D13" = "0.03070"
    "D14" = "0.09358"
    "D15" = "3.728"
    "D16" = "0.001606"
    "D17" = "0.04777"
    "D18" = "0.0005791"
    "D19" = "0.006377"
    "D20" = "0.005047"
    "D21" = "0.001029"
    "D22" = "0.0001501"
    "D23" = "3.706"
    "D24" = "2.178"
    "D25" = "0.3305"
    "D26" = "0.1308"
    "D28" = "0.0003001"
    "D40" = "0.04027"
    "B41" = "BKEXToken"
    "C41" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "D41" = "0.1072"
    "E41" = "40BKEXTokenBKK"
    "B42" = "CEJI"
    "C42" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
    "D42" = "0.002712"
    "E42" = "41CEJICEJI"
    "B43" = "KickToken"
    "C43" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "D43" = "0.002943"
    "E43" = "42KickTokenKICKWorstin24h"
    "D44" = "0.008020"
    "D45" = "0.00005804"
    "D47" = "0.5001"
    "E47" = "46CoinbaseStockTokenCOIN"
    "D49" = "0.00002100"
    "D50" = "0.01010"
}

foreach ($ref in $newValues.Keys) {
    $value = $newValues[$ref]
    $cell = $ws.Range($ref)

    $looksNumeric = $value -match '^[0-9]+(\.[0-9]+)?$'

    if ($looksNumeric) {
        # Force text storage, write the literal digits, then strip the
        # temporary formatting so the cell ends up with no style index
        # (same as its original, un-styled state).
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.ClearFormats()
    } else {
        # Not a numeric-looking string (coin name / URL / label) - a plain
        # assignment already keeps it as text.
        $cell.Value = $value
    }
}
